# [Mod] Reaction Data Cell 수정
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReactionDatas")

# Header row
$ws.Range("D1").Value = "honor"

# Data rows: ActionType, Script, Weight, Reward
$data = @(
    @("Perfecct", "아..! 살 것 같다..", 1.5, 5),
    @("Perfecct", "감사합니다! 정말 맛있어 보이네요!", 1.5, 5),
    @("Perfecct", "둘이 먹다 하나 쪄 죽어도 모를 맛!", 1.5, 5),
    @("Perfecct", "와! 이게 아이스크림이구나!", 1.5, 5),
    @("Perfecct", "최고에요! ", 1.5, 5),
    @("Good", "감사합니다.", 1, 2),
    @("Good", "뭔가 다른 것 같은데.. 그래도 잘 먹을게요.", 1, 2),
    @("Good", "잘 먹겠습니다.", 1, 2),
    @("Good", "음.. 이런 맛이구나.", 1, 2),
    @("Good", "맛있네요.", 1, 2),
    @("Bad", "이게.. 맞나요..?", 0.5, -5),
    @("Bad", "에잉 이게 아닌데..쯧..", 0.5, -5),
    @("Bad", "뭐야.. 장사할 생각 없어요?", 0.5, -5),
    @("Bad", "내 아이스크림이 아니야!", 0.5, -5),
    @("Bad", "이럴거면 주문 왜 받은거죠?", 0.5, -5)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}

$ws.Columns.Item(2).ColumnWidth = 35.75

$ws.Range("G14").Select()
